$d = $word.ActiveDocument

# --- 1) "...utilizzando JavaFX che permette..." : merge the split JavaFX run
#        (removes the spell-check proofErr wrap around JavaFX) ---
$d.Content.Find.Execute("utilizzando JavaFX che permette", $false, $true, $false, $false, $false, `
    $true, 1, $false, "utilizzando JavaFX che permette", 2) | Out-Null

# --- 2) "Essa memorizza le carte, partite e tornei su filesystem..." :
#        insert " ad eliminazione diretta" right after "tornei" ---
$rng = $d.Content
$rng.Find.Execute("carte, partite e tornei", $false, $true, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$found = $rng.Duplicate
$found.Collapse(0)
$found.InsertAfter(" ad eliminazione diretta")

# --- 3) "Passa alla directory contenente il file .jar con il comando cd " :
#        merge the split ".jar" run (removes proofErr wrap around "jar") ---
$d.Content.Find.Execute("alla directory contenente il file .jar con il comando cd", $false, $true, $false, $false, $false, `
    $true, 1, $false, "alla directory contenente il file .jar con il comando cd", 2) | Out-Null

# --- 4) "Esegui il file gara utilizzando il comando java -jar " :
#        merge the split "jar" run ---
$d.Content.Find.Execute("Esegui il file gara utilizzando il comando java -jar ", $false, $true, $false, $false, $false, `
    $true, 1, $false, "Esegui il file gara utilizzando il comando java -jar ", 2) | Out-Null

# --- 5) "...il seguente comando mvn clean javafx:run." :
#        merge mvn / clean / javafx:run runs (removes proofErr + gramStart/End wraps) ---
$d.Content.Find.Execute("mvn clean javafx:run", $false, $true, $false, $false, $false, `
    $true, 1, $false, "mvn clean javafx:run", 2) | Out-Null

# --- 6) "...repository pubblico di Github" : merge the split "Github" run ---
$d.Content.Find.Execute("repository pubblico di Github", $false, $true, $false, $false, $false, `
    $true, 1, $false, "repository pubblico di Github", 2) | Out-Null

# --- 7) Add justified alignment to the last (empty) paragraph ---
$last = $d.Paragraphs.Last
$last.Range.ParagraphFormat.Alignment = 3
